# Update Work Week and Social Spending
# (Actually: refresh the GDP per Capita data series for Turkmenistan with
#  a newer data vintage, extending coverage from 2010 through 2016.)
#
# Column E on the "Data" sheet stores its numbers as *text* (shared
# strings), matching the source export. To reproduce that without Excel
# silently reinterpreting the string as a number (which would also leave
# a stray cell style behind), we temporarily mark the cell as Text,
# assign the value, and then clear the formatting back to the default
# so the workbook's style table is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

function Set-TextValue {
    param($Range, $Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# --- Existing rows whose "Data" value changes -----------------------------

Set-TextValue $ws.Range("E2") "7693"

Set-TextValue $ws.Range("E9")  "6763"
Set-TextValue $ws.Range("E10") "6666"
Set-TextValue $ws.Range("E11") "6575"
Set-TextValue $ws.Range("E12") "6535"
Set-TextValue $ws.Range("E13") "6314"
Set-TextValue $ws.Range("E14") "6186"
Set-TextValue $ws.Range("E15") "6204"
Set-TextValue $ws.Range("E16") "6194"
Set-TextValue $ws.Range("E17") "6601"
Set-TextValue $ws.Range("E18") "5915"
Set-TextValue $ws.Range("E19") "5796"
Set-TextValue $ws.Range("E20") "5457.75692573531"
Set-TextValue $ws.Range("E21") "4587.76365541108"
Set-TextValue $ws.Range("E22") "4604.38415341451"
Set-TextValue $ws.Range("E23") "3769.62476924548"
Set-TextValue $ws.Range("E24") "3473.68018927826"
Set-TextValue $ws.Range("E25") "3226.98483638665"
Set-TextValue $ws.Range("E26") "2856.14860685559"
Set-TextValue $ws.Range("E27") "3037.70433134861"
Set-TextValue $ws.Range("E28") "3527.16046040025"
Set-TextValue $ws.Range("E29") "4173.19459953499"
Set-TextValue $ws.Range("E30") "5014.06245100545"
Set-TextValue $ws.Range("E31") "5795.6409686451"
Set-TextValue $ws.Range("E32") "6784.56100135347"
Set-TextValue $ws.Range("E33") "7781.74382479198"
Set-TextValue $ws.Range("E34") "8791.67584575014"
Set-TextValue $ws.Range("E35") "9752.35547849855"
Set-TextValue $ws.Range("E36") "10829.3757801065"
Set-TextValue $ws.Range("E37") "12426.718788804"
Set-TextValue $ws.Range("E38") "13188.4478549552"
Set-TextValue $ws.Range("E39") "14397.1490969721"

# --- New rows: 2011-2016 ----------------------------------------------------

$newRows = @(
    @{ Row = 40; Year = 2011; Value = "16518" },
    @{ Row = 41; Year = 2012; Value = "18135" },
    @{ Row = 42; Year = 2013; Value = "19751" },
    @{ Row = 43; Year = 2014; Value = "21546" },
    @{ Row = 44; Year = 2015; Value = "22675" },
    @{ Row = 45; Year = 2016; Value = "23813" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 795
    $ws.Range("B$row").Value = "Turkmenistan"
    $ws.Range("C$row").Value = "GDP per Capita"
    $ws.Range("D$row").Value = $r.Year
    Set-TextValue $ws.Range("E$row") $r.Value
}
